# Rename sheets to align CM module case-study workbook with the main
# data model / new get_data() naming scheme.
#
#   FreshwaterSources              -> ExternalWaterSources
#   FreshwaterSourcingAvailability -> ExtWaterSourcingAvailability
#   FreshSourcingCost              -> ExternalSourcingCost
#   BeneficialReuseRevenue         -> BeneficialReuseCredit

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("FreshwaterSources").Name = "ExternalWaterSources"
$wb.Worksheets.Item("FreshwaterSourcingAvailability").Name = "ExtWaterSourcingAvailability"
$wb.Worksheets.Item("FreshSourcingCost").Name = "ExternalSourcingCost"
$wb.Worksheets.Item("BeneficialReuseRevenue").Name = "BeneficialReuseCredit"
